# Auto-generated edit script applying numeric corrections to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H129").Value = 1466.0444
$ws.Range("J129").Value = 1485.7954
$ws.Range("L129").Value = 4457.3862
$ws.Range("N129").Value = -14457.3862

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16942.117
$ws.Range("I32").Value = 17900.807
$ws.Range("J32").Value = 7035.6665
$ws.Range("K32").Value = 17900.807
$ws.Range("L32").Value = 7035.6665
$ws.Range("M32").Value = -17613.807
$ws.Range("N32").Value = -7609.6665
$ws.Range("H88").Value = 103488.9
$ws.Range("I88").Value = 1522.6
$ws.Range("J88").Value = 205455.2
$ws.Range("K88").Value = 1522.6
$ws.Range("L88").Value = 205455.2
$ws.Range("M88").Value = -1116.6
$ws.Range("N88").Value = -206267.2
$ws.Range("H91").Value = 103488.9
$ws.Range("I91").Value = 1522.6
$ws.Range("J91").Value = 205455.2
$ws.Range("K91").Value = 1522.6
$ws.Range("L91").Value = 205455.2
$ws.Range("M91").Value = -118.5999999999999
$ws.Range("N91").Value = -208263.2
$ws.Range("H102").Value = 1265.3889
$ws.Range("I102").Value = 1055.7142
$ws.Range("J102").Value = 1999.25
$ws.Range("K102").Value = 1055.7142
$ws.Range("L102").Value = 1999.25
$ws.Range("M102").Value = 566.2858000000001
$ws.Range("N102").Value = -5243.25
$ws.Range("H122").Value = 2456.111
$ws.Range("I122").Value = 2469.375
$ws.Range("J122").Value = 2350
$ws.Range("K122").Value = 7408.125
$ws.Range("L122").Value = 7050
$ws.Range("M122").Value = -4958.125
$ws.Range("N122").Value = -11950

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 231
$ws.Range("I11").Value = 231
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 231
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -91
$ws.Range("N11").Value = ""
$ws.Range("H20").Value = 3040
$ws.Range("I20").Value = 2805.3845
$ws.Range("K20").Value = 2805.3845
$ws.Range("M20").Value = -2558.3845
$ws.Range("H86").Value = 1587.3877
$ws.Range("I86").Value = 1366.079
$ws.Range("J86").Value = 2351.9092
$ws.Range("K86").Value = 1366.079
$ws.Range("L86").Value = 2351.9092
$ws.Range("M86").Value = -243.079
$ws.Range("N86").Value = -4597.9092
$ws.Range("H89").Value = 1587.3877
$ws.Range("I89").Value = 1366.079
$ws.Range("J89").Value = 2351.9092
$ws.Range("K89").Value = 6830.395
$ws.Range("L89").Value = 11759.546
$ws.Range("M89").Value = -1214.395
$ws.Range("N89").Value = -22991.546
$ws.Range("H94").Value = 2363.875
$ws.Range("I94").Value = 1078.1666
$ws.Range("K94").Value = 1078.1666
$ws.Range("M94").Value = -627.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 998.6875
$ws.Range("I122").Value = 813.8333
$ws.Range("J122").Value = 1553.25
$ws.Range("K122").Value = 2441.4999
$ws.Range("L122").Value = 4659.75
$ws.Range("M122").Value = 8.500100000000202
$ws.Range("N122").Value = -9559.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 100
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = ""
$ws.Range("H131").Value = 785.14
$ws.Range("J131").Value = 788.051
$ws.Range("L131").Value = 2364.153
$ws.Range("N131").Value = -12444.153

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 15025501
$ws.Range("I3").Value = 50000000
$ws.Range("J3").Value = 3367334.8
$ws.Range("K3").Value = 50000000
$ws.Range("L3").Value = 3367334.8
$ws.Range("M3").Value = -49999884
$ws.Range("N3").Value = -3367566.8
$ws.Range("H42").Value = 25600
$ws.Range("J42").Value = 25600
$ws.Range("L42").Value = 25600
$ws.Range("N42").Value = -26570
$ws.Range("H80").Value = 3339.818
$ws.Range("I80").Value = 2696.5
$ws.Range("J80").Value = 3875.9167
$ws.Range("K80").Value = 2696.5
$ws.Range("L80").Value = 3875.9167
$ws.Range("M80").Value = -1698.5
$ws.Range("N80").Value = -5871.9167
$ws.Range("H83").Value = 3339.818
$ws.Range("I83").Value = 2696.5
$ws.Range("J83").Value = 3875.9167
$ws.Range("K83").Value = 13482.5
$ws.Range("L83").Value = 19379.5835
$ws.Range("M83").Value = -8490.5
$ws.Range("N83").Value = -29363.5835
$ws.Range("H97").Value = 2380.7778
$ws.Range("I97").Value = 966.7273
$ws.Range("K97").Value = 966.7273
$ws.Range("M97").Value = -470.7273
$ws.Range("H113").Value = 1915.0625
$ws.Range("I113").Value = 1588.4615
$ws.Range("J113").Value = 2138.5264
$ws.Range("K113").Value = 1588.4615
$ws.Range("L113").Value = 2138.5264
$ws.Range("M113").Value = 581.5385000000001
$ws.Range("N113").Value = -6478.526400000001
$ws.Range("H115").Value = 25600
$ws.Range("J115").Value = 25600
$ws.Range("L115").Value = 25600
$ws.Range("N115").Value = -27950
$ws.Range("H122").Value = 2412.2354
$ws.Range("I122").Value = 1788.375
$ws.Range("K122").Value = 5365.125
$ws.Range("M122").Value = -2915.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6157
$ws.Range("I22").Value = 5475.25
$ws.Range("J22").Value = 7066
$ws.Range("K22").Value = 5475.25
$ws.Range("L22").Value = 7066
$ws.Range("M22").Value = -5180.25
$ws.Range("N22").Value = -7656
$ws.Range("H27").Value = 6157
$ws.Range("I27").Value = 5475.25
$ws.Range("J27").Value = 7066
$ws.Range("K27").Value = 5475.25
$ws.Range("L27").Value = 7066
$ws.Range("M27").Value = -5368.25
$ws.Range("N27").Value = -7280
$ws.Range("H61").Value = 3637.1538
$ws.Range("I61").Value = 1424.375
$ws.Range("J61").Value = 7177.6
$ws.Range("K61").Value = 1424.375
$ws.Range("L61").Value = 7177.6
$ws.Range("M61").Value = -1222.375
$ws.Range("N61").Value = -7581.6
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").Value = ""
$ws.Range("H113").Value = 3637.1538
$ws.Range("I113").Value = 1424.375
$ws.Range("J113").Value = 7177.6
$ws.Range("K113").Value = 1424.375
$ws.Range("L113").Value = 7177.6
$ws.Range("M113").Value = 745.625
$ws.Range("N113").Value = -11517.6
$ws.Range("H122").Value = 1511235.5
$ws.Range("I122").Value = 1963035.6
$ws.Range("K122").Value = 5889106.800000001
$ws.Range("M122").Value = -5886656.800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 50000000
$ws.Range("I8").Value = 50000000
$ws.Range("K8").Value = 50000000
$ws.Range("M8").Value = -49999860
$ws.Range("H21").Value = 575.5
$ws.Range("J21").Value = 1680
$ws.Range("L21").Value = 1680
$ws.Range("N21").Value = -2150
$ws.Range("H35").Value = 575.5
$ws.Range("J35").Value = 1680
$ws.Range("L35").Value = 1680
$ws.Range("N35").Value = -2260
$ws.Range("H81").Value = 83334530
$ws.Range("I81").Value = 1299.7273
$ws.Range("J81").Value = 1000000000
$ws.Range("K81").Value = 2599.4546
$ws.Range("L81").Value = 2000000000
$ws.Range("M81").Value = -1538.4546
$ws.Range("N81").Value = -2000002122
$ws.Range("H84").Value = 83334530
$ws.Range("I84").Value = 1299.7273
$ws.Range("J84").Value = 1000000000
$ws.Range("K84").Value = 12997.273
$ws.Range("L84").Value = 10000000000
$ws.Range("M84").Value = -7693.273000000001
$ws.Range("N84").Value = -10000010608
$ws.Range("H122").Value = 1694.5714
$ws.Range("I122").Value = 1704.9524
$ws.Range("J122").Value = 1663.4286
$ws.Range("K122").Value = 5114.857199999999
$ws.Range("L122").Value = 4990.2858
$ws.Range("M122").Value = -2664.857199999999
$ws.Range("N122").Value = -9890.2858
